# Fruta / hortaliza, semanal
# Insert a new weekly price record for Ajo (Terminal Hortofrutícola Agro
# Chillán) above the existing row 123, shifting the remaining records
# (old rows 123-146) down by one row (new rows 124-147).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 123:146 down to 124:147, making room for the new record.
$ws.Rows("123:123").Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A123").Value = 7
$ws.Range("B123").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C123").Value = "Ñuble"
$ws.Range("D123").Value = 44511
$ws.Range("E123").Value = 16
$ws.Range("F123").Value = 100112003
$ws.Range("G123").Value = "Ajo"
$ws.Range("H123").Value = "Chino"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 100
$ws.Range("K123").Value = 15000
$ws.Range("L123").Value = 16000
$ws.Range("M123").Value = 15500
$ws.Range("N123").Value = "$/caja 10 kilos"
$ws.Range("O123").Value = "China"
$ws.Range("P123").Value = 1550
$ws.Range("Q123").Value = 10
$ws.Range("R123").Value = "Hortaliza"
